# Order creation scenario along with code refactoring
# qualis_users_credentials.xlsx - "test" sheet update:
#  - opsadmin password changed (and its stray hyperlink removed)
#  - "forgot password user" row password value refreshed, extra custom
#    formatting dropped (refactor)
#  - new "fund manager" row added
#  - new "invalid user" row added (with a hyperlink on the email cell)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# All hyperlinks on this sheet are rebuilt from scratch below, so clear the
# existing collection first (per-hyperlink deletion isn't reliable here).
$ws.Hyperlinks.Delete()

# Drop ALL pre-existing cell formatting first (the old blue/underlined
# Hyperlink style on rows 2-5 and the one-off "JetBrains Mono" / centered
# style that used to live on row 6 only), so nothing extra gets layered on
# top while we rebuild the hyperlinks below.
$ws.Range("A1:C6").Style = "Normal"

# ---- Row 1: headers (unchanged) ----
$ws.Range("A1").Value = "User"
$ws.Range("B1").Value = "Email"
$ws.Range("C1").Value = "Password"

# ---- Row 2: superuser (unchanged) ----
$ws.Range("A2").Value = "superuser"
$ws.Range("B2").Value = "autoqualissuperuser@praemium.com"
$ws.Range("C2").Value = "QS@superuserPSS123!"

# ---- Row 3: opsadmin (password replaced, no more hyperlink on password) ----
$ws.Range("A3").Value = "opsadmin"
$ws.Range("B3").Value = "autoqualisuser_opsadmin@praemium.com"
$ws.Range("C3").Value = "!Qualis1!"

# ---- Row 4: viewuser (unchanged) ----
$ws.Range("A4").Value = "viewuser"
$ws.Range("B4").Value = "autoqualisuser_view@praemium.com"
$ws.Range("C4").Value = "QS@viewuserPSS123!"

# ---- Row 5: finadviser (unchanged) ----
$ws.Range("A5").Value = "finadviser"
$ws.Range("B5").Value = "autofinadviser@praemium.com"
$ws.Range("C5").Value = "QS@FinAdviser1PSS123!"

# ---- Row 6: forgot password user (new throwaway password, plain formatting) ----
$ws.Range("A6").Value = "forgot password user"
$ws.Range("B6").Value = "qataskdemoaccnt@gmail.com"
$ws.Range("C6").Value = "newPT_638*235"

# ---- Row 7: fund manager (new row) ----
$ws.Range("A7").Value = "fund manager"
$ws.Range("B7").Value = "autofundmanager@praemium.com"
$ws.Range("C7").Value = "QS@fundManagerPSS123!"

# ---- Row 8: invalid user (new row) ----
$ws.Range("A8").Value = "invalid user"
$ws.Range("B8").Value = "test@test.com"
$ws.Range("C8").Value = "pass12345!."

# Re-create the mailto hyperlinks that still belong on the sheet. Cell text
# is already set above, so Add() just wires up the link without clobbering it.
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:autoqualissuperuser@praemium.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:QS@superuserPSS123!") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:autoqualisuser_opsadmin@praemium.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:autoqualisuser_view@praemium.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:QS@viewuserPSS123!") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B6"), "mailto:qataskdemoaccnt@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B8"), "mailto:test@test.com") | Out-Null

# Hyperlinks.Add stamps its own "link" look-alike style on every cell it
# touches; strip that back off (plain text, no special formatting) except
# for B8 which should keep the normal blue/underlined Hyperlink style.
$ws.Range("B2").Style = "Normal"
$ws.Range("C2").Style = "Normal"
$ws.Range("B3").Style = "Normal"
$ws.Range("B4").Style = "Normal"
$ws.Range("C4").Style = "Normal"
$ws.Range("B6").Style = "Normal"
$ws.Range("B8").Style = "Hyperlink"

# Match the saved selection state (cursor resting on the newly added last row).
$ws.Range("A8").Select()
